$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.465.81"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "2.510.77"
$ws.Range("E3").Value = "  -4.96%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").Value = "2.511.59"
$ws.Range("E9").Value = "  -4.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.53%  "
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "2.968.84"
$ws.Range("E14").Value = "  -5.06%  "
$ws.Range("D15").Value = "69.403.19"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").Value = "2.505.88"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("E23").Value = "  -5.93%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.38%  "
$ws.Range("D28").Value = "2.630.52"
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "0.0₃0897"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "457.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.72%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.79%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.48%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("E50").Value = "  -4.93%  "
$ws.Range("E51").Value = "  -1.71%  "
